$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.321.49'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '2.275.10'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.40%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.98%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.68'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Value = '2.627.78'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('D16').Value = '2.272.75'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.793'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.98%  '
$ws.Range('D18').Value = '42.234.10'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +3.66%  '
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.95'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('E39').Value = '  +2.90%  '
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.10'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.68%  '
$ws.Range('E42').Value = '  +13.93%  '
$ws.Range('D43').Value = '1.998.19'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0285'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.98'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '52.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '91.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.04%  '
